# Apply "end of loop" progress update to the Bill Summary sheet:
# randomized Qty-executed-upto-date (column C) values for rows 8-17,
# plus the recomputed Upto date Amount (column G) for the rows that
# carry a rate, and the recalculated Grand Total / Net Payable cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C - "Qty executed upto date"
$ws.Range("C8").Value  = 66
$ws.Range("C9").Value  = 85
$ws.Range("C10").Value = 58
$ws.Range("C11").Value = 88
$ws.Range("C12").Value = 98
$ws.Range("C13").Value = 96
$ws.Range("C14").Value = 39
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 100
$ws.Range("C17").Value = 12

# Column G - "Upto date Amount" (Qty upto date * Rate), stored as formatted text
$ws.Range("G9").Value  = "21760.00"
$ws.Range("G10").Value = "27376.00"
$ws.Range("G11").Value = "58256.00"
$ws.Range("G13").Value = "13056.00"
$ws.Range("G14").Value = "897.00"

# Grand total / net payable rows recalculated from the new amounts
$ws.Range("G19").Value = "121345.00"
$ws.Range("H19").Value = "121345.00"
$ws.Range("G21").Value = "121345.00"
$ws.Range("H21").Value = "121345.00"
